# Fruta / hortaliza, semanal
# Insert two new weekly rows of data at the top of the "1033" data block,
# shifting the existing rows (1033-1127) down by two (to 1035-1129), then
# populate the two newly-inserted rows (1033, 1034) with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 1033; existing data shifts down.
$ws.Rows("1033:1034").Insert()

# New row 1033 - Calidad "Primera"
$ws.Range("A1033").Value = 3
$ws.Range("B1033").Value = "Femacal de La Calera"
$ws.Range("C1033").Value = "Coquimbo"
$ws.Range("D1033").Value = 45132
$ws.Range("E1033").Value = 5
$ws.Range("F1033").Value = 100112008
$ws.Range("G1033").Value = "Coliflor"
$ws.Range("H1033").Value = "Sin especificar"
$ws.Range("I1033").Value = "Primera"
$ws.Range("J1033").Value = 3500
$ws.Range("K1033").Value = 700
$ws.Range("L1033").Value = 750
$ws.Range("M1033").Value = 723
$ws.Range("N1033").Value = "$/unidad"
$ws.Range("O1033").Value = "Provincia de Quillota"
$ws.Range("P1033").Value = 723
$ws.Range("Q1033").Value = 1
$ws.Range("R1033").Value = "Hortaliza"

# New row 1034 - Calidad "Segunda"
$ws.Range("A1034").Value = 3
$ws.Range("B1034").Value = "Femacal de La Calera"
$ws.Range("C1034").Value = "Coquimbo"
$ws.Range("D1034").Value = 45132
$ws.Range("E1034").Value = 5
$ws.Range("F1034").Value = 100112008
$ws.Range("G1034").Value = "Coliflor"
$ws.Range("H1034").Value = "Sin especificar"
$ws.Range("I1034").Value = "Segunda"
$ws.Range("J1034").Value = 500
$ws.Range("K1034").Value = 550
$ws.Range("L1034").Value = 550
$ws.Range("M1034").Value = 550
$ws.Range("N1034").Value = "$/unidad"
$ws.Range("O1034").Value = "Provincia de Quillota"
$ws.Range("P1034").Value = 550
$ws.Range("Q1034").Value = 1
$ws.Range("R1034").Value = "Hortaliza"
